# Gradebook update: record Chapter 1 quiz/test grades and add the
# "Chapter 2" assignments/grades block (copy of the Chapter 1 layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Chapter 1: fill in the quiz1 and test1 grades that were missing ---
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 101

# Re-enter the final-grade formula so it re-evaluates now that its
# precedents (H9/J9) are no longer errors.
$ws.Range("J11").Formula = "=F9*0.4+H9*0.25+J9*0.3+5"

# --- Chapter 2 block (rows 13-22), mirroring rows 2-11 ---

# Section title
$ws.Range("F13:K13").Merge()
$ws.Range("F13").Value = "Assignments and Grades Received - Chapter 2"

# Column group headers: Homework / Quiz / Test
$ws.Range("F14:G14").Merge()
$ws.Range("F14").Value = "Homework"
$ws.Range("H14:I14").Merge()
$ws.Range("H14").Value = "Quiz"
$ws.Range("J14:K14").Merge()
$ws.Range("J14").Value = "Test"

# Sub headers: Assignment / Grade
$ws.Range("F15").Value = "Assignment"
$ws.Range("G15").Value = "Grade"
$ws.Range("H15").Value = "Assignment"
$ws.Range("I15").Value = "Grade"
$ws.Range("J15").Value = "Assignment"
$ws.Range("K15").Value = "Grade"

# Assignment rows
$ws.Range("F16").Value = "hw1"
$ws.Range("G16").Value = 94
$ws.Range("H16").Value = "quiz1"
$ws.Range("J16").Value = "test1"
$ws.Range("F17").Value = "hw2"

# Overall grade labels
$ws.Range("F19:G19").Merge()
$ws.Range("F19").Value = "Overall Homework Grade"
$ws.Range("H19:I19").Merge()
$ws.Range("H19").Value = "Overall Quiz Grade"
$ws.Range("J19:K19").Merge()
$ws.Range("J19").Value = "Overall Test Grade"

# Overall grade formulas
$ws.Range("F20:G20").Merge()
$ws.Range("F20").Formula = "=AVERAGE(G16:G17)"
$ws.Range("H20:I20").Merge()
$ws.Range("H20").Formula = "=AVERAGE(I16:I17)"
$ws.Range("J20:K20").Merge()
$ws.Range("J20").Formula = "=AVERAGE(K16:K17)"

# Final chapter grade
$ws.Range("F22:I22").Merge()
$ws.Range("F22").Value = "Chapter 1 Grade"
$ws.Range("J22:K22").Merge()
$ws.Range("J22").Formula = "=F20*0.4+H20*0.25+J20*0.3+5"

# Match the author's last-saved selection
$ws.Range("L8").Select()
